$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2043.6111
$ws.Range("I19").Value = 4018.75
$ws.Range("J19").Value = 463.5
$ws.Range("K19").Value = 4018.75
$ws.Range("L19").Value = 463.5
$ws.Range("M19").Value = -3843.75
$ws.Range("N19").Value = -813.5

$ws.Range("H98").Value = 1704.5526
$ws.Range("J98").Value = 3796.6365
$ws.Range("L98").Value = 3796.6365
$ws.Range("N98").Value = -6792.636500000001

$ws.Range("H113").Value = 2836.1365
$ws.Range("I113").Value = 2635.25
$ws.Range("J113").Value = 2950.9285
$ws.Range("K113").Value = 2635.25
$ws.Range("L113").Value = 2950.9285
$ws.Range("M113").Value = 618.75
$ws.Range("N113").Value = -9458.9285

$ws.Range("H122").Value = 1704.5526
$ws.Range("J122").Value = 3796.6365
$ws.Range("L122").Value = 11389.9095
$ws.Range("N122").Value = -16289.9095

$ws.Range("H125").Value = 2166.125
$ws.Range("J125").Value = 4994.6
$ws.Range("L125").Value = 44951.4
$ws.Range("N125").Value = -49871.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 12500
$ws.Range("J24").Value = 12500
$ws.Range("L24").Value = 12500
$ws.Range("N24").Value = -13248

$ws.Range("H28").Value = 24958.846
$ws.Range("I28").Value = 8155
$ws.Range("K28").Value = 8155
$ws.Range("M28").Value = -7963

$ws.Range("H31").Value = 10377.625
$ws.Range("I31").Value = 4288.7144
$ws.Range("K31").Value = 4288.7144
$ws.Range("M31").Value = -3994.7144

$ws.Range("H32").Value = 3468.5386
$ws.Range("I32").Value = 3490.5403
$ws.Range("K32").Value = 3490.5403
$ws.Range("M32").Value = -3203.5403

$ws.Range("H74").Value = 1681.5883
$ws.Range("I74").Value = 1681.5883
$ws.Range("K74").Value = 1681.5883
$ws.Range("M74").Value = -807.5882999999999

$ws.Range("H77").Value = 1681.5883
$ws.Range("I77").Value = 1681.5883
$ws.Range("K77").Value = 8407.941499999999
$ws.Range("M77").Value = -4039.941499999999

$ws.Range("H82").Value = 25000
$ws.Range("J82").Value = 25000
$ws.Range("L82").Value = 25000
$ws.Range("N82").Value = -25722

$ws.Range("H85").Value = 25000
$ws.Range("J85").Value = 25000
$ws.Range("L85").Value = 25000
$ws.Range("N85").Value = -27496

$ws.Range("H86").Value = 43749.25
$ws.Range("J86").Value = 43749.25
$ws.Range("L86").Value = 43749.25
$ws.Range("N86").Value = -46121.25

$ws.Range("H89").Value = 43749.25
$ws.Range("J89").Value = 43749.25
$ws.Range("L89").Value = 131247.75
$ws.Range("N89").Value = -143103.75

$ws.Range("H93").Value = 27500
$ws.Range("I93").Value = 15000
$ws.Range("K93").Value = 15000
$ws.Range("M93").Value = -12504

$ws.Range("H96").Value = 29856.857
$ws.Range("J96").Value = 29856.857
$ws.Range("L96").Value = 29856.857
$ws.Range("N96").Value = -35348.857

$ws.Range("H97").Value = 47668220
$ws.Range("I97").Value = 50001630
$ws.Range("J97").Value = 1000000
$ws.Range("K97").Value = 50001630
$ws.Range("L97").Value = 1000000
$ws.Range("M97").Value = -50001134
$ws.Range("N97").Value = -1000992

$ws.Range("H99").Value = 24958.846
$ws.Range("I99").Value = 8155
$ws.Range("K99").Value = 8155
$ws.Range("M99").Value = -5160

$ws.Range("H100").Value = 12500
$ws.Range("J100").Value = 12500
$ws.Range("L100").Value = 12500
$ws.Range("N100").Value = -14664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 30
$ws.Range("I4").Value = 30
$ws.Range("K4").Value = 30
$ws.Range("M4").Value = 82

$ws.Range("H31").Value = 2919.4285
$ws.Range("I31").Value = 2512.25
$ws.Range("J31").Value = 3462.3333
$ws.Range("K31").Value = 2512.25
$ws.Range("L31").Value = 3462.3333
$ws.Range("M31").Value = -2217.25
$ws.Range("N31").Value = -4052.3333

$ws.Range("H34").Value = 2919.4285
$ws.Range("I34").Value = 2512.25
$ws.Range("J34").Value = 3462.3333
$ws.Range("K34").Value = 2512.25
$ws.Range("L34").Value = 3462.3333
$ws.Range("M34").Value = -2310.25
$ws.Range("N34").Value = -3866.3333

$ws.Range("H132").Value = 252119.64
$ws.Range("I132").Value = 8683.532999999999
$ws.Range("J132").Value = 916036.25
$ws.Range("K132").Value = 26050.599
$ws.Range("L132").Value = 2748108.75
$ws.Range("M132").Value = -23520.599
$ws.Range("N132").Value = -2753168.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 2240.111
$ws.Range("I22").Value = 1032.6
$ws.Range("K22").Value = 3097.8
$ws.Range("M22").Value = -2928.8

$ws.Range("H27").Value = 2240.111
$ws.Range("I27").Value = 1032.6
$ws.Range("K27").Value = 3097.8
$ws.Range("M27").Value = -2995.8

$ws.Range("H41").Value = 871
$ws.Range("I41").Value = 324.44446
$ws.Range("J41").Value = 2100.75
$ws.Range("K41").Value = 973.33338
$ws.Range("L41").Value = 6302.25
$ws.Range("M41").Value = -635.33338
$ws.Range("N41").Value = -6978.25

$ws.Range("H44").Value = 1184.5
$ws.Range("I44").Value = 252.66667
$ws.Range("K44").Value = 758.00001
$ws.Range("M44").Value = -360.00001

$ws.Range("H69").Value = 4618.3
$ws.Range("I69").Value = 4111.857
$ws.Range("K69").Value = 12335.571
$ws.Range("M69").Value = -11524.571

$ws.Range("H72").Value = 4618.3
$ws.Range("I72").Value = 4111.857
$ws.Range("K72").Value = 37006.713
$ws.Range("M72").Value = -32950.713

$ws.Range("H92").Value = 707.7857
$ws.Range("I92").Value = 790.25
$ws.Range("J92").Value = 674.8
$ws.Range("K92").Value = 2370.75
$ws.Range("L92").Value = 2024.4
$ws.Range("M92").Value = -1122.75
$ws.Range("N92").Value = -4520.4

$ws.Range("H97").Value = 778.6429000000001
$ws.Range("I97").Value = 668.2857
$ws.Range("K97").Value = 2004.8571
$ws.Range("M97").Value = -1508.8571

$ws.Range("H117").Value = 1451.0834
$ws.Range("J117").Value = 1583.6666
$ws.Range("L117").Value = 4750.9998
$ws.Range("N117").Value = -11634.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 17766.846
$ws.Range("I132").Value = 19998.883
$ws.Range("J132").Value = 2589
$ws.Range("K132").Value = 59996.649
$ws.Range("L132").Value = 7767
$ws.Range("M132").Value = -57466.649
$ws.Range("N132").Value = -12827

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 28332
$ws.Range("J75").Value = 29999
$ws.Range("L75").Value = 29999
$ws.Range("N75").Value = -31871

$ws.Range("H78").Value = 28332
$ws.Range("J78").Value = 29999
$ws.Range("L78").Value = 89997
$ws.Range("N78").Value = -99357

$ws.Range("H81").Value = 2935325.2
$ws.Range("I81").Value = 3638888.8
$ws.Range("J81").Value = 3810.6667
$ws.Range("K81").Value = 7277777.6
$ws.Range("L81").Value = 7621.3334
$ws.Range("M81").Value = -7276716.6
$ws.Range("N81").Value = -9743.3334

$ws.Range("H84").Value = 2935325.2
$ws.Range("I84").Value = 3638888.8
$ws.Range("J84").Value = 3810.6667
$ws.Range("K84").Value = 36388888
$ws.Range("L84").Value = 38106.667
$ws.Range("M84").Value = -36383584
$ws.Range("N84").Value = -48714.667

$ws.Range("H139").Value = 68023
$ws.Range("J139").Value = 68023
$ws.Range("L139").Value = 68023
$ws.Range("N139").Value = -78303
